# Actualización automática 2025-11-12 17:30:08
# Applies the updated sales figures to the three report sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (sales per client / product group)
# ---------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("C4").Value  = 518.4
$wsGrupo.Range("M4").Value  = 1554.43

$wsGrupo.Range("L5").Value  = 835.22

$wsGrupo.Range("D29").Value = 475.2
$wsGrupo.Range("I29").Value = 313.2
$wsGrupo.Range("K29").Value = 812.16
$wsGrupo.Range("L29").Value = 1677.67

$wsGrupo.Range("D37").Value = 457.92
$wsGrupo.Range("M37").Value = 1960.7

$wsGrupo.Range("H47").Value = 782.1

# Row 56 holds the "N de 54" summary counters per product column.
$wsGrupo.Range("C56").Value = "2 de 54"
$wsGrupo.Range("D56").Value = "5 de 54"
$wsGrupo.Range("H56").Value = "1 de 54"
$wsGrupo.Range("I56").Value = "2 de 54"
$wsGrupo.Range("L56").Value = "5 de 54"
$wsGrupo.Range("M56").Value = "11 de 54"

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL" (sales per client / month)
# ---------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F4").Value  = 3130.56
$wsMensual.Range("F5").Value  = 984.91
$wsMensual.Range("F29").Value = 3278.23
$wsMensual.Range("F37").Value = 2418.62
$wsMensual.Range("F47").Value = 782.1
$wsMensual.Range("F60").Value = 20736.88

# ---------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" (budget compliance per product group)
# ---------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column F (CUMPLIMIENTO) width narrowed slightly.
$wsCumpl.Columns.Item(6).ColumnWidth = 23.1

$wsCumpl.Range("D2").Value  = 907.2
$wsCumpl.Range("E2").Value  = 1919.46
$wsCumpl.Range("F2").Value  = 0.3209441531701726

$wsCumpl.Range("D3").Value  = 2799.36
$wsCumpl.Range("E3").Value  = 3823.9
$wsCumpl.Range("F3").Value  = 0.4226559126472462

$wsCumpl.Range("D6").Value  = 782.1
$wsCumpl.Range("E6").Value  = 2125.48368146026
$wsCumpl.Range("F6").Value  = 0.2689862393254354

$wsCumpl.Range("D7").Value  = 336.6
$wsCumpl.Range("E7").Value  = 983.4
$wsCumpl.Range("F7").Value  = 0.255

$wsCumpl.Range("D10").Value = 811.65
$wsCumpl.Range("E10").Value = 3500.35
$wsCumpl.Range("F10").Value = 0.1882305194805195

$wsCumpl.Range("D11").Value = 6050.76
$wsCumpl.Range("E11").Value = 8185.23
$wsCumpl.Range("F11").Value = 0.4250326110091395

$wsCumpl.Range("D12").Value = 9025.809999999999
$wsCumpl.Range("E12").Value = 55918.19
$wsCumpl.Range("F12").Value = 0.1389783505789603

$wsCumpl.Range("D14").Value = 20713.48
$wsCumpl.Range("E14").Value = 78242.77685923838
$wsCumpl.Range("F14").Value = 0.2093195585344761
